# The six species-observation records in rows 2-7 of the Artfynd sheet get
# reshuffled onto different rows (the underlying records are unchanged, only
# which row they occupy changes). Concretely, using (new row) <- (old row):
#   2 <- 6   3 <- 5   4 <- 7   5 <- 4   6 <- 3   7 <- 2
#
# Rather than moving whole rows around (which would also disturb untouched
# cells), only the cells whose value actually differs between the old and
# new layout are written here, matching the published change set exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# The "Antal" column (I) holds numbers written as text (e.g. "256"); Excel
# would otherwise auto-convert such digit-only strings to real numbers, so
# force a text format on those cells before writing them.
function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Row 2 (now holds what used to be row 6's record) ---
Set-Cell 2 1 111545323          # A2 Id
Set-TextCell 2 9 "2"            # I2 Antal
Set-Cell 2 16 "Orsa Viborg, intill en grupp med hålträd, Dlr"  # P2 Lokalnamn
Set-Cell 2 17 480427.8053356989 # Q2 Ost
Set-Cell 2 18 6772811.198980245 # R2 Nord

# --- Row 3 (now holds what used to be row 5's record) ---
Set-Cell 3 1 111545414          # A3 Id
Set-TextCell 3 9 "9"            # I3 Antal
Set-Cell 3 10 "registreringar"  # J3 Enhet
Set-Cell 3 16 "Orsa Viborg, glänta i mitten av skogsparti, Dlr" # P3 Lokalnamn
Set-Cell 3 17 480487.2503558649 # Q3 Ost
Set-Cell 3 18 6772784.264016891 # R3 Nord

# --- Row 4 (now holds what used to be row 7's record) ---
Set-Cell 4 1 111545401          # A4 Id
Set-Cell 4 2 57487               # B4 Taxonsorteringsordning
Set-Cell 4 4 "NT"                # D4 Rödlistade
Set-Cell 4 5 205998               # E4 TaxonId
Set-Cell 4 6 "Nordfladdermus"    # F4 Artnamn
Set-Cell 4 7 "Eptesicus nilssonii" # G4 Vetenskapligt namn
Set-Cell 4 8 "(A.Keyserling & Blasius, 1839)" # H4 Auktor
Set-TextCell 4 9 "6"             # I4 Antal
Set-Cell 4 16 "Orsa Viborg, glänta i mitten av skogsparti, Dlr" # P4 Lokalnamn
Set-Cell 4 17 480487.2503558649  # Q4 Ost
Set-Cell 4 18 6772784.264016891  # R4 Nord

# --- Row 5 (now holds what used to be row 4's record) ---
Set-Cell 5 1 111543957            # A5 Id
Set-TextCell 5 9 "1"              # I5 Antal
Set-Cell 5 16 "Orsa Viborg, glänta i skogsparti, Dlr" # P5 Lokalnamn
Set-Cell 5 17 480406.6045043401   # Q5 Ost
Set-Cell 5 18 6772745.04339793    # R5 Nord

# --- Row 6 (now holds what used to be row 3's record) ---
Set-Cell 6 1 111545328             # A6 Id
Set-Cell 6 2 57494                 # B6 Taxonsorteringsordning
Set-Cell 6 4 "LC"                  # D6 Rödlistade
Set-Cell 6 5 205992                # E6 TaxonId
Set-Cell 6 6 "Vattenfladdermus"    # F6 Artnamn
Set-Cell 6 7 "Myotis daubentonii"  # G6 Vetenskapligt namn
Set-Cell 6 8 "(Kuhl, 1817)"        # H6 Auktor
Set-TextCell 6 9 "1"               # I6 Antal

# --- Row 7 (now holds what used to be row 2's record) ---
Set-Cell 7 1 111543968             # A7 Id
Set-TextCell 7 9 "256"             # I7 Antal
$ws.Cells.Item(7, 10).Value = ""   # J7 Enhet (cleared)
Set-Cell 7 16 "Orsa Viborg, glänta i skogsparti, Dlr" # P7 Lokalnamn
Set-Cell 7 17 480406.6045043401    # Q7 Ost
Set-Cell 7 18 6772745.04339793     # R7 Nord
